# BOM_JLCSMT_DaisySeedGuitarPedal1590b-SMD.xlsx update
# JLCPCB now stocks CPC1018N regularly, so the BOM is updated from the
# previously special-ordered CPC1019N part to CPC1018N, along with its
# new JLCPCB part number (C2760117 -> C133069).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 holds the relay opto-isolator part: Comment (A), Designator (B),
# Footprint (C), JLCPCB Part # (D).
$ws.Range("A27").Value = "CPC1018N"
$ws.Range("D27").Value = "C133069"

# The active selection in the saved file moved from B28 to A28.
$ws.Range("A28").Select()
